# Add cantrals by cantons
# Restructure Sheet1: merge the old two-row header (labels row + units row)
# into a single header row, add new idx/idx2/Name/Date columns info and
# rename the units so every data row moves up by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The old sheet had:
#   row 1 = top labels (Alimentation / pompes) / Hiver / Eté / Année) spread
#           across E1,G1,I1,J1,K1
#   row 2 = unit labels ( (m3/s) (MW) (MW) (GWh) (GWh) (GWh) ) in F2:K2
#   rows 3-14 = the 12 data rows
#
# The new sheet has a single header row (row 1) with explicit column
# names, and the 12 data rows shifted up to rows 2-13 unchanged.
#
# Remove the old units row (row 2); this shifts the 12 data rows up by one,
# so they land on rows 2-13 exactly like the target layout.
$ws.Rows("2:2").Delete()

# Rewrite row 1 completely with the new header text.
$ws.Cells.Item(1, 1).Value = "idx"
$ws.Cells.Item(1, 2).Value = "idx2"
$ws.Cells.Item(1, 3).Value = "Name"
$ws.Cells.Item(1, 4).Value = "Date Start"
$ws.Cells.Item(1, 5).Value = "Date End"
$ws.Cells.Item(1, 6).Value = "(m3/s)"
$ws.Cells.Item(1, 7).Value = "(MW1)"
$ws.Cells.Item(1, 8).Value = "(MW2)"
$ws.Cells.Item(1, 9).Value = "(GWh) Winter"
$ws.Cells.Item(1, 10).Value = "(GWh) Summer"
$ws.Cells.Item(1, 11).Value = "(GWh) Year"

# Columns A-E of the header carry no special styling in the target file.
$ws.Cells.Item(1, 1).ClearFormats()
$ws.Cells.Item(1, 2).ClearFormats()
$ws.Cells.Item(1, 3).ClearFormats()
$ws.Cells.Item(1, 4).ClearFormats()
$ws.Cells.Item(1, 5).ClearFormats()

# Columns F-K of the header keep the workbook's "label" font (Arial 9).
for ($col = 6; $col -le 11; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.ClearFormats()
    $cell.Font.Name = "Arial"
    $cell.Font.Size = 9
}

# Match the new selection left behind by the edit (first data row).
$ws.Range("A2:K2").Select()

Write-Output "Restructured Sheet1 header/data rows"
